# Update "想去人数" (F column) figures across the 展览, 本地生活 and 全部类型
# sheets to reflect the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4, 6).Value = 3773
$ws1.Cells.Item(6, 6).Value = 5319
$ws1.Cells.Item(7, 6).Value = 602
$ws1.Cells.Item(8, 6).Value = 439
$ws1.Cells.Item(9, 6).Value = 233
$ws1.Cells.Item(10, 6).Value = 1076
$ws1.Cells.Item(12, 6).Value = 147
$ws1.Cells.Item(14, 6).Value = 737
$ws1.Cells.Item(15, 6).Value = 371
$ws1.Cells.Item(18, 6).Value = 183
$ws1.Cells.Item(21, 6).Value = 6106
$ws1.Cells.Item(25, 6).Value = 7056
$ws1.Cells.Item(29, 6).Value = 376
$ws1.Cells.Item(30, 6).Value = 761
$ws1.Cells.Item(32, 6).Value = 326
$ws1.Cells.Item(33, 6).Value = 141
$ws1.Cells.Item(35, 6).Value = 1177
$ws1.Cells.Item(39, 6).Value = 935
$ws1.Cells.Item(40, 6).Value = 1166

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 1157

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value = 1157
$ws4.Cells.Item(7, 6).Value = 3773
$ws4.Cells.Item(9, 6).Value = 5319
$ws4.Cells.Item(10, 6).Value = 602
$ws4.Cells.Item(11, 6).Value = 439
$ws4.Cells.Item(12, 6).Value = 233
$ws4.Cells.Item(13, 6).Value = 1076
$ws4.Cells.Item(15, 6).Value = 147
$ws4.Cells.Item(17, 6).Value = 737
$ws4.Cells.Item(18, 6).Value = 371
$ws4.Cells.Item(22, 6).Value = 183
$ws4.Cells.Item(25, 6).Value = 6106
$ws4.Cells.Item(29, 6).Value = 7056
$ws4.Cells.Item(33, 6).Value = 376
$ws4.Cells.Item(34, 6).Value = 761
$ws4.Cells.Item(36, 6).Value = 326
$ws4.Cells.Item(38, 6).Value = 141
$ws4.Cells.Item(40, 6).Value = 1177
$ws4.Cells.Item(44, 6).Value = 935
$ws4.Cells.Item(45, 6).Value = 1166
